# The workbook has 10 sheets; the last one is "Лист1" (sheetId=1), which is
# the currently active/selected tab. The edit duplicates that sheet,
# inserting the copy immediately before it (Excel names the duplicate
# "Лист1 (11)"). The new sheet becomes the active tab, while the original
# "Лист1" sheet keeps its place as the very last tab but is no longer
# selected. The new sheet's cell values are then updated to a new training
# sample (a 20x20 grid of 0/1 values, similar to but not identical to the
# other "Лист1 (N)" samples already in the workbook).

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

# Duplicate the last sheet, inserting the new copy right before it.
# Excel auto-names it "Лист1 (11)" and, per observed COM semantics here,
# the newly inserted copy becomes the active/selected sheet/tab.
$lastSheet.Copy($lastSheet)

# The new copy is now the sheet immediately before the original last sheet.
$newSheet = $wb.Worksheets.Item($sheetCount)

# Target data for the new sample sheet (20 rows x 20 cols, A1:T20).
$rowsData = @(
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,1,1,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,1,1,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,1,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,1,1,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
)

for ($r = 1; $r -le 20; $r++) {
    $rowVals = $rowsData[$r - 1]
    for ($c = 1; $c -le 20; $c++) {
        $newSheet.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Make sure the new sheet is the active selection/tab.
$newSheet.Activate()
$newSheet.Select()
